$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 142; this shifts old rows 142:257 down to 143:258
$ws.Rows.Item(142).Insert()

# Fill in the newly inserted row 142 with the new weekly record
$ws.Range("A142").Value = 5
$ws.Range("B142").Value = "Macroferia Regional de Talca"
$ws.Range("C142").Value = "Maule"
$ws.Range("D142").Value = 44729
$ws.Range("E142").Value = 7
$ws.Range("F142").Value = 100112045
$ws.Range("G142").Value = "Zapallo"
$ws.Range("H142").Value = "Camote"
$ws.Range("I142").Value = "1a (guarda)"
$ws.Range("J142").Value = 800
$ws.Range("K142").Value = 400
$ws.Range("L142").Value = 400
$ws.Range("M142").Value = 400
$ws.Range("N142").Value = '$/kilo (volumen en unidades)'
$ws.Range("O142").Value = "Región del Maule"
$ws.Range("P142").Value = 400
$ws.Range("Q142").Value = 1
$ws.Range("R142").Value = "Hortaliza"
